$wb = $excel.ActiveWorkbook

# --- Rename header labels on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Copy header style (bold, bordered, centered) from the Weekly Quantity header row
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Copy the date-formatted style from the Weekly Quantity date column
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A34").PasteSpecial(-4122)

# Header row values
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Data values
$dsVals = @(44934.99999999999,44955.99999999999,44962.99999999999,44976.99999999999,44983.99999999999,44997.99999999999,45004.99999999999,45011.99999999999,45018.99999999999,45025.99999999999,45032.99999999999,45039.99999999999,45046.99999999999,45053.99999999999,45060.99999999999,45074.99999999999,45088.99999999999,45095.99999999999,45102.99999999999,45116.99999999999,45123.99999999999,45130.99999999999,45137.99999999999,45151.99999999999,45158.99999999999,45165.99999999999,45172.99999999999,45179.99999999999,45186.99999999999,45193.99999999999,45200.99999999999,45207.99999999999,45214.99999999999)
$poVals = @(189,180,178,172,169,163,160,157,154,151,148,145,142,139,137,131,125,122,119,113,110,107,104,98,96,93,90,87,84,81,78,75,72)
$lowerVals = @(60.05047004731592,58.4870503785903,59.51537437510949,41.3946678587992,39.39901090786323,37.27359425536577,32.496288248044,30.98504833201339,36.91970631909386,26.837341592928,31.47985220687629,31.32792496270732,27.14289328172596,-0.6309359632685171,6.976751104264886,6.657702842916057,0.2483766575570169,-7.18524203954901,-4.654192741056016,-4.020719733805517,-7.598237994079819,-22.01086369808385,-15.03854194033911,-18.41423338702263,-21.71695737039369,-22.7525794055598,-37.35191630477956,-35.49433307756857,-31.77990119933714,-38.27664496214371,-39.03997032739916,-47.55756785707407,-43.40630521889556)
$upperVals = @(320.7070097236362,306.2308363186708,295.5692808801745,297.5592920962633,294.281082405762,287.717614685934,283.0328730360682,278.0002230806307,270.3878437324735,271.1229989800068,257.7332821738713,276.4390234639757,265.1539613400683,266.5652163494067,253.8257288409463,246.6752207795767,244.1486003133364,241.4652814474676,229.8846095774234,232.8035862816917,236.6937499258818,227.7393469286588,226.2715853147117,227.0846786849954,209.4819992051784,217.6790358187251,222.4576039408016,209.228264630636,210.6935759524032,210.8625199828471,200.8699286350494,199.0174897106794,203.9040616000127)

for ($i = 0; $i -lt $dsVals.Length; $i++) {
    $row = $i + 2
    $wsForecast.Cells.Item($row, 1).Value = $dsVals[$i]
    $wsForecast.Cells.Item($row, 2).Value = $poVals[$i]
    $wsForecast.Cells.Item($row, 3).Value = $lowerVals[$i]
    $wsForecast.Cells.Item($row, 4).Value = $upperVals[$i]
}

$wsForecast.Range("A1").Select() | Out-Null
